$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) values must remain text exactly as scraped (e.g. "1.000" not 1).
# Force text storage via NumberFormat, then strip the style back to Normal so no
# stray style index is left on the cell (matches original unstyled inlineStr cells).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "24.044.35"
$ws.Range("E2").Value = "  +0.62%  "
Set-TextValue $ws.Range("D3") "1.660.47"
$ws.Range("E3").Value = "  +2.18%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.57%  "
Set-TextValue $ws.Range("D5") "309.82"
$ws.Range("E5").Value = "  +1.03%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  -0.50%  "
Set-TextValue $ws.Range("D7") "0.3908"
$ws.Range("E7").Value = "  +0.04%  "
Set-TextValue $ws.Range("D8") "0.3871"
$ws.Range("E8").Value = "  +1.57%  "
Set-TextValue $ws.Range("D9") "51.25"
$ws.Range("E9").Value = "  +3.19%  "
Set-TextValue $ws.Range("D10") "1.367"
$ws.Range("E10").Value = "  +0.71%  "
Set-TextValue $ws.Range("D11") "1.000"
$ws.Range("E11").Value = "  -0.58%  "
Set-TextValue $ws.Range("D12") "0.08512"
$ws.Range("E12").Value = "  +0.95%  "
Set-TextValue $ws.Range("D13") "23.95"
$ws.Range("E13").Value = "  +1.04%  "
Set-TextValue $ws.Range("D14") "7.220"
$ws.Range("E14").Value = "  +3.36%  "
Set-TextValue $ws.Range("D15") "8.040"
$ws.Range("E15").Value = "  +8.39%  "
Set-TextValue $ws.Range("D16") "0.00001315"
$ws.Range("E16").Value = "  +3.62%  "
Set-TextValue $ws.Range("D17") "1.659.14"
$ws.Range("E17").Value = "  +1.81%  "
Set-TextValue $ws.Range("D18") "94.64"
$ws.Range("E18").Value = "  +2.14%  "
Set-TextValue $ws.Range("D19") "0.06997"
$ws.Range("E19").Value = "  +1.28%  "
Set-TextValue $ws.Range("D20") "19.97"
$ws.Range("E20").Value = "  +0.84%  "
Set-TextValue $ws.Range("D21") "6.977"
$ws.Range("E21").Value = "  +2.09%  "
Set-TextValue $ws.Range("D22") "1.001"
$ws.Range("E22").Value = "  -0.58%  "
Set-TextValue $ws.Range("D23") "13.67"
$ws.Range("E23").Value = "  +2.45%  "
Set-TextValue $ws.Range("D24") "24.040.44"
$ws.Range("E24").Value = "  +0.54%  "
Set-TextValue $ws.Range("D25") "2.483"
$ws.Range("E25").Value = "  +3.46%  "
Set-TextValue $ws.Range("D26") "3.091"
$ws.Range("E26").Value = "  +7.84%  "
Set-TextValue $ws.Range("D27") "22.24"
$ws.Range("E27").Value = "  +0.68%  "
Set-TextValue $ws.Range("D28") "153.86"
$ws.Range("E28").Value = "  -2.37%  "
Set-TextValue $ws.Range("D29") "139.99"
$ws.Range("E29").Value = "  +1.04%  "
Set-TextValue $ws.Range("D30") "5.318"
$ws.Range("E30").Value = "  +1.19%  "
Set-TextValue $ws.Range("D31") "7.979"
$ws.Range("E31").Value = "  +4.32%  "
Set-TextValue $ws.Range("D32") "2.488"
$ws.Range("E32").Value = "  +1.08%  "
Set-TextValue $ws.Range("D33") "1.841.77"
$ws.Range("E33").Value = "  +2.22%  "
Set-TextValue $ws.Range("D34") "1.047"
$ws.Range("E34").Value = "  +9.24%  "
Set-TextValue $ws.Range("D35") "0.08145"
$ws.Range("E35").Value = "  +2.79%  "
Set-TextValue $ws.Range("D36") "0.03006"
$ws.Range("E36").Value = "  +4.45%  "
Set-TextValue $ws.Range("D37") "11.16"
$ws.Range("E37").Value = "  +9.31%  "
Set-TextValue $ws.Range("D38") "6.749"
$ws.Range("E38").Value = "  +2.49%  "
Set-TextValue $ws.Range("D39") "0.2708"
$ws.Range("E39").Value = "  +2.13%  "
Set-TextValue $ws.Range("D40") "0.09164"
$ws.Range("E40").Value = "  +0.36%  "
Set-TextValue $ws.Range("D41") "13.71"
$ws.Range("E41").Value = "  +2.52%  "
Set-TextValue $ws.Range("D42") "0.7573"
$ws.Range("E42").Value = "  +2.03%  "
Set-TextValue $ws.Range("D43") "1.423"
$ws.Range("E43").Value = "  +0.39%  "
Set-TextValue $ws.Range("D44") "16.43"
$ws.Range("E44").Value = "  +3.25%  "
Set-TextValue $ws.Range("D45") "0.7035"
$ws.Range("E45").Value = "  +3.25%  "
Set-TextValue $ws.Range("D46") "2.499"
$ws.Range("E46").Value = "  +2.69%  "
Set-TextValue $ws.Range("D47") "4.097"
$ws.Range("E47").Value = "  +0.64%  "
Set-TextValue $ws.Range("D48") "0.9997"
$ws.Range("E48").Value = "  -0.62%  "
Set-TextValue $ws.Range("D49") "0.08300"
$ws.Range("E49").Value = "  +0.73%  "
Set-TextValue $ws.Range("D50") "135.83"
$ws.Range("E50").Value = "  +2.50%  "
Set-TextValue $ws.Range("D51") "1.241"
$ws.Range("E51").Value = "  -0.69%  "
